$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's entire column range (1-16384) carries a "Text" number format
# style, so a bare Range.Value numeric assignment gets coerced to a text
# string. Reset the target cells to the "Normal" style first so the new
# values are stored as real numbers, matching the rest of the numeric data.

# New column D header (" Oct 06")
$ws.Range("D1").Style = "Normal"
$ws.Range("D1").Value = " Oct 06"

# New column D values for existing rows 2-4
$ws.Range("D2").Style = "Normal"
$ws.Range("D2").Value = 30

$ws.Range("D3").Style = "Normal"
$ws.Range("D3").Value = 5

$ws.Range("D4").Style = "Normal"
$ws.Range("D4").Value = 5

# New rows 5-7 with label in column A and value 0 in column D
$ws.Range("A5").Style = "Normal"
$ws.Range("A5").Value = "Deleted Trailers"
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").Value = 0

$ws.Range("A6").Style = "Normal"
$ws.Range("A6").Value = "Deleted Trucks"
$ws.Range("D6").Style = "Normal"
$ws.Range("D6").Value = 0

$ws.Range("A7").Style = "Normal"
$ws.Range("A7").Value = "Trailer"
$ws.Range("D7").Style = "Normal"
$ws.Range("D7").Value = 0
